$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "StarExp" column header is renamed to "SoulStone" as part of the
# Cookie Enhance Lv / Star API work.
$ws.Range("C1").Value = "SoulStone"

# Reflect the author's final selection in the saved sheet view.
$ws.Range("C2").Select()
